$d = $word.ActiveDocument

# Find the "Docente(s) Responsável(eis)" heading paragraph.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Docente(s) Responsável(eis)", $false, $false,
                                    $false, $false, $false, $true, 1, $false,
                                    "", 0)

if ($found) {
    $target = $searchRange.Paragraphs(1)
} else {
    # Fallback: scan paragraphs directly if Find didn't collapse onto one.
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith("Docente(s) Responsável(eis)")) {
            $target = $p
            break
        }
    }
}

# Create a brand new paragraph right after it, then fill it in (via raw
# OOXML so we get the exact ListBullet-styled run/line-break structure)
# with the instructor roster.
$null = $target.Range.InsertParagraphAfter()
$newPara = $target.Next()

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
       '<w:r><w:t>7043088 - Ana Karine Furtado de Carvalho</w:t><w:br/></w:r>' +
       '<w:r><w:t>7926291 - Célia Regina Tomachuk dos Santos Catuogno</w:t><w:br/></w:r>' +
       '<w:r><w:t>4893449 - Débora Souza Alvim</w:t><w:br/></w:r>' +
       '<w:r><w:t>8855158 - Morun Bernardino Neto</w:t><w:br/></w:r>' +
       '<w:r><w:t>7455355 - Robson da Silva Rocha</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

$null = $newPara.Range.InsertXML($xml)
